# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Irak overtakes Belgica/Ecuador in the ranking (rows 31-33 swap country
#    labels while carrying fresh case numbers for the new #1 of the trio)
#  - Refreshed totals for several other countries (USA, India, Argentina,
#    Portugal, Serbia, Uzbekistan, Rep. of Macedonia, Islandia, Hong Kong,
#    Niger)
#  - The "Datos actualizados" timestamp footer moves from 14:33 to 15:50

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Estados Unidos (row 4) ----
$ws.Range("B4").Value = 2985897
$ws.Range("C4").Value = 2969
$ws.Range("D4").Value = 1289836
$ws.Range("E4").Value = 1563451

# ---- India (row 6) ----
$ws.Range("B6").Value = 704607
$ws.Range("C6").Value = 6771
$ws.Range("D6").Value = 429907
$ws.Range("E6").Value = 254918
$ws.Range("G6").Value = 82
$ws.Range("H6").Value = 19782

# ---- Argentina (row 26) ----
$ws.Range("D26").Value = 28531
$ws.Range("E26").Value = 47761
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 1523

# ---- Rows 31-33: Irak jumps ahead of Belgica and Ecuador ----
# Row 31 now shows Irak with new, higher case counts.
$ws.Range("A31").Value = "Irak"
$ws.Range("B31").Value = 62275
$ws.Range("C31").Value = 1796
$ws.Range("D31").Value = 34741
$ws.Range("E31").Value = 24967
$ws.Range("G31").Value = 94
$ws.Range("H31").Value = 2567

# Row 32 now shows Belgica, carrying what used to be row 31's numbers.
$ws.Range("A32").Value = "Belgica"
$ws.Range("B32").Value = 62016
$ws.Range("C32").Value = 107
$ws.Range("D32").Value = 17091
$ws.Range("E32").Value = 35154
$ws.Range("H32").Value = 9771

# Row 33 now shows Ecuador, carrying what used to be row 32's numbers.
$ws.Range("A33").Value = "Ecuador"
$ws.Range("B33").Value = 61958
$ws.Range("D33").Value = 28722
$ws.Range("E33").Value = 28455
$ws.Range("H33").Value = 4781

# ---- Portugal (row 42) ----
$ws.Range("B42").Value = 44129
$ws.Range("C42").Value = 232
$ws.Range("D42").Value = 29166
$ws.Range("E42").Value = 13343
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 1620

# ---- Serbia (row 62) ----
$ws.Range("B62").Value = 16420
$ws.Range("C62").Value = 289
$ws.Range("D62").Value = 13366
$ws.Range("E62").Value = 2737
$ws.Range("G62").Value = 6
$ws.Range("H62").Value = 317

# ---- Uzbekistan (row 71) ----
$ws.Range("B71").Value = 10284
$ws.Range("C71").Value = 264
$ws.Range("E71").Value = 3664
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 36

# ---- Republica de Macedonia (row 83) ----
$ws.Range("B83").Value = 7124
$ws.Range("C83").Value = 78
$ws.Range("D83").Value = 3199
$ws.Range("E83").Value = 3579
$ws.Range("G83").Value = 5
$ws.Range("H83").Value = 346

# ---- Islandia (row 116) ----
$ws.Range("B116").Value = 1866
$ws.Range("C116").Value = 3
$ws.Range("D116").Value = 1840
$ws.Range("E116").Value = 16

# ---- Hong Kong (row 127) ----
$ws.Range("B127").Value = 1286
$ws.Range("C127").Value = 17
$ws.Range("D127").Value = 1157
$ws.Range("E127").Value = 122

# ---- Niger (row 134) ----
$ws.Range("B134").Value = 1093
$ws.Range("C134").Value = 5
$ws.Range("D134").Value = 968
$ws.Range("E134").Value = 57

# ---- Footer timestamp ----
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 15:50"
